{"js": "// Office.js (Word JavaScript API) edit matching the target diff.\n//\n// Original body: a single paragraph\n//   \"OBJETIVO GERAL \" + <br/> + <br/> + <br/> + \"Realizar o levantamento...\"\n//\n// Target body: the same first paragraph trimmed down to\n//   \"OBJETIVO GERAL \" + <br/> + <br/> + \"Objetivo Geral\"\n// followed by a series of new paragraphs holding the \"Objetivo Geral\" /\n// \"Objetivos Espec\u00edficos\" write-up.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\n\n// Wipe the existing run content of the first paragraph (keeps the\n// paragraph mark / paragraph-level formatting), then rebuild it with the\n// new wording: title line, two manual line breaks, \"Objetivo Geral\".\nconst firstRange = firstParagraph.getRange();\nfirstRange.insertText(\"\", \"Replace\");\nawait context.sync();\n\nfirstParagraph.insertText(\"OBJETIVO GERAL \", \"End\");\nawait context.sync();\nfirstParagraph.insertBreak(\"Line\", \"End\");\nawait context.sync();\nfirstParagraph.insertBreak(\"Line\", \"End\");\nawait context.sync();\nfirstParagraph.insertText(\"Objetivo Geral\", \"End\");\nawait context.sync();\n\n// New paragraphs inserted after the (now shortened) first paragraph.\nconst newParagraphTexts = [\n  \"Desenvolver um sistema de gerenciamento de chamados, denominado UpDesk, com integra\u00e7\u00e3o de intelig\u00eancia artificial, a fim de otimizar o processo de triagem, atendimento e resolu\u00e7\u00e3o de demandas t\u00e9cnicas, promovendo maior agilidade, organiza\u00e7\u00e3o e efici\u00eancia no suporte ao usu\u00e1rio. \",\n  \" Objetivos Espec\u00edficos\",\n  \"Levantar e documentar os requisitos funcionais e n\u00e3o funcionais do sistema com base em um cen\u00e1rio real.\",\n  \"Modelar o sistema utilizando diagramas UML (casos de uso, classes, sequ\u00eancia e implanta\u00e7\u00e3o), aplicando boas pr\u00e1ticas de engenharia de software.\",\n  \"Implementar o sistema utilizando tecnologias adequadas, garantindo seguran\u00e7a, escalabilidade e controle de acesso.\",\n  \"Integrar um m\u00f3dulo de intelig\u00eancia artificial para sugerir solu\u00e7\u00f5es automaticamente durante a abertura de chamados.\",\n  \"Criar um prot\u00f3tipo visual no Figma para valida\u00e7\u00e3o da interface com os usu\u00e1rios.\",\n  \"Aplicar princ\u00edpios da Lei Geral de Prote\u00e7\u00e3o de Dados (LGPD) para assegurar a privacidade e seguran\u00e7a das informa\u00e7\u00f5es dos usu\u00e1rios.\",\n  \"Desenvolver e aplicar roteiros de teste e scripts de dados para validar o funcionamento do sistema em ambiente de homologa\u00e7\u00e3o.IV).\",\n];\n\nlet lastParagraph = firstParagraph;\nfor (const text of newParagraphTexts) {\n  lastParagraph = lastParagraph.insertParagraph(text, \"After\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop edit matching the target diff.\n#\n# Original body: a single paragraph\n#   \"OBJETIVO GERAL \" + <br/> + <br/> + <br/> + \"Realizar o levantamento...\"\n#\n# Target body: the same first paragraph trimmed down to\n#   \"OBJETIVO GERAL \" + <br/> + <br/> + \"Objetivo Geral\"\n# followed by a series of new paragraphs holding the \"Objetivo Geral\" /\n# \"Objetivos Espec\u00edficos\" write-up.\n\n$d = $word.ActiveDocument\n$p1 = $d.Paragraphs.Item(1)\n\n# 1) Remove one of the three manual line breaks right after \"OBJETIVO GERAL \".\n#    (\"OBJETIVO GERAL \" is 15 characters, so the breaks start at offset 15\n#    relative to the start of the paragraph/document.)\n$full = $p1.Range\n$breakStart = $full.Start + 15\n$breakRange = $d.Range($breakStart, $breakStart + 1)\n$breakRange.Delete()\n\n# 2) Replace the long \"Realizar o levantamento...\" sentence with the new\n#    short title text \"Objetivo Geral\".\n$oldText = \"Realizar o levantamento e an\u00e1lise de requisitos para a cria\u00e7\u00e3o de um sistema de suporte t\u00e9cnico inteligente, que utilize Intelig\u00eancia Artificial para triagem inicial, categoriza\u00e7\u00e3o autom\u00e1tica de chamados e sugest\u00e3o de solu\u00e7\u00f5es, visando otimizar o tempo de resposta e reduzir a sobrecarga da equipe de TI. A LGPD deve ser aplicada a todos os dados pessoais tratados no sistema. O desenvolvimento ocorrer\u00e1 no pr\u00f3ximo semestre (PIM IV).\"\n$targetRange = $p1.Range\n$null = $targetRange.Find.Execute($oldText)\n$targetRange.Text = \"Objetivo Geral\"\n\n# 3) Append the new paragraphs (Objetivo Geral write-up + Objetivos\n#    Espec\u00edficos list) right after the first paragraph.\n$newParagraphTexts = @(\n  \"Desenvolver um sistema de gerenciamento de chamados, denominado UpDesk, com integra\u00e7\u00e3o de intelig\u00eancia artificial, a fim de otimizar o processo de triagem, atendimento e resolu\u00e7\u00e3o de demandas t\u00e9cnicas, promovendo maior agilidade, organiza\u00e7\u00e3o e efici\u00eancia no suporte ao usu\u00e1rio. \",\n  \" Objetivos Espec\u00edficos\",\n  \"Levantar e documentar os requisitos funcionais e n\u00e3o funcionais do sistema com base em um cen\u00e1rio real.\",\n  \"Modelar o sistema utilizando diagramas UML (casos de uso, classes, sequ\u00eancia e implanta\u00e7\u00e3o), aplicando boas pr\u00e1ticas de engenharia de software.\",\n  \"Implementar o sistema utilizando tecnologias adequadas, garantindo seguran\u00e7a, escalabilidade e controle de acesso.\",\n  \"Integrar um m\u00f3dulo de intelig\u00eancia artificial para sugerir solu\u00e7\u00f5es automaticamente durante a abertura de chamados.\",\n  \"Criar um prot\u00f3tipo visual no Figma para valida\u00e7\u00e3o da interface com os usu\u00e1rios.\",\n  \"Aplicar princ\u00edpios da Lei Geral de Prote\u00e7\u00e3o de Dados (LGPD) para assegurar a privacidade e seguran\u00e7a das informa\u00e7\u00f5es dos usu\u00e1rios.\",\n  \"Desenvolver e aplicar roteiros de teste e scripts de dados para validar o funcionamento do sistema em ambiente de homologa\u00e7\u00e3o.IV).\"\n)\n\n$cur = $p1\nforeach ($t in $newParagraphTexts) {\n  $cur.Range.InsertParagraphAfter()\n  $cur = $d.Paragraphs.Item($cur.Index + 1)\n  $cur.Range.Text = $t\n}\n"}
